$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 705; existing rows 705:796 shift down to 706:797.
$ws.Rows.Item(705).Insert()

# Populate the newly inserted row 705 with the new data point.
$ws.Range("A705").Value = 6
$ws.Range("B705").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C705").Value = "Metropolitana"
$ws.Range("D705").Value = 45154
$ws.Range("E705").Value = 13
$ws.Range("F705").Value = 100112044
$ws.Range("G705").Value = "Perejil"
$ws.Range("H705").Value = "Sin especificar"
$ws.Range("I705").Value = "Primera"
$ws.Range("J705").Value = 320
$ws.Range("K705").Value = 11000
$ws.Range("L705").Value = 12000
$ws.Range("M705").Value = 11469
$ws.Range("N705").Value = "`$/docena de atados"
$ws.Range("O705").Value = "Región Metropolitana"
$ws.Range("P705").Value = 3823
$ws.Range("Q705").Value = 3
$ws.Range("R705").Value = "Hortaliza"
